$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / field names) - H1 and I1 swap order (SamplePortion now before Result)
$ws.Range("A1").Value = "Operator"
$ws.Range("B1").Value = "SampleID"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "LaboratoryOperatingMode"
$ws.Range("E1").Value = "CriticalApparatusCriticalSoftware"
$ws.Range("F1").Value = "CriticalProduct"
$ws.Range("G1").Value = "RawDataPathway"
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"
$ws.Range("J1").Value = "CookingTime"
$ws.Range("K1").Value = "Temperature"

# Row 2 (type row) - add units to float/integer temperature types
$ws.Range("A2").Value = "#string"
$ws.Range("B2").Value = "#string"
$ws.Range("C2").Value = "#date"
$ws.Range("D2").Value = "#string"
$ws.Range("E2").Value = "#string"
$ws.Range("F2").Value = "#string"
$ws.Range("G2").Value = "#string"
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"
$ws.Range("J2").Value = "#integer"
$ws.Range("K2").Value = "#integer,  unit:celsius"

# Row 3 (new descriptions / enum row)
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#TempsCuisson"
$ws.Range("K3").Value = "#Temperature"
